$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3123.6538
$ws.Range("I64").Value = 2967.1667
$ws.Range("K64").Value = 2967.1667
$ws.Range("M64").Value = -2719.1667

$ws.Range("H67").Value = 3123.6538
$ws.Range("I67").Value = 2967.1667
$ws.Range("K67").Value = 2967.1667
$ws.Range("M67").Value = -2109.1667

$ws.Range("H76").Value = 3706879
$ws.Range("I76").Value = 7410268
$ws.Range("J76").Value = 3490
$ws.Range("K76").Value = 7410268
$ws.Range("L76").Value = 3490
$ws.Range("M76").Value = -7409953
$ws.Range("N76").Value = -4120

$ws.Range("H79").Value = 3706879
$ws.Range("I79").Value = 7410268
$ws.Range("J79").Value = 3490
$ws.Range("K79").Value = 7410268
$ws.Range("L79").Value = 3490
$ws.Range("M79").Value = -7409176
$ws.Range("N79").Value = -5674

$ws.Range("H115").Value = 875.36365
$ws.Range("J115").Value = 1159
$ws.Range("L115").Value = 3477
$ws.Range("N115").Value = -6611

$ws.Range("H118").Value = 2089.1667
$ws.Range("I118").Value = 270
$ws.Range("J118").Value = 2254.5454
$ws.Range("K118").Value = 810
$ws.Range("L118").Value = 6763.6362
$ws.Range("M118").Value = 847
$ws.Range("N118").Value = -10077.6362

$ws.Range("H123").Value = 46296
$ws.Range("J123").Value = 46296
$ws.Range("L123").Value = 46296
$ws.Range("N123").Value = -56096

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3575
$ws.Range("I63").Value = 1620
$ws.Range("J63").Value = 6833.3335
$ws.Range("K63").Value = 1620
$ws.Range("L63").Value = 6833.3335
$ws.Range("M63").Value = -934
$ws.Range("N63").Value = -8205.333500000001

$ws.Range("H66").Value = 3575
$ws.Range("I66").Value = 1620
$ws.Range("J66").Value = 6833.3335
$ws.Range("K66").Value = 8100
$ws.Range("L66").Value = 34166.6675
$ws.Range("M66").Value = -4668
$ws.Range("N66").Value = -41030.6675

$ws.Range("H102").Value = 3141.2083
$ws.Range("I102").Value = 1831.8
$ws.Range("K102").Value = 1831.8
$ws.Range("M102").Value = -209.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H19").Value = 2000
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H35").Value = 11100
$ws.Range("J35").Value = 11100
$ws.Range("L35").Value = 11100
$ws.Range("N35").Value = -11720

$ws.Range("H82").Value = 15076.714
$ws.Range("I82").Value = 4889.25
$ws.Range("J82").Value = 28660
$ws.Range("K82").Value = 4889.25
$ws.Range("L82").Value = 28660
$ws.Range("M82").Value = -4506.25
$ws.Range("N82").Value = -29426

$ws.Range("H85").Value = 15076.714
$ws.Range("I85").Value = 4889.25
$ws.Range("J85").Value = 28660
$ws.Range("K85").Value = 4889.25
$ws.Range("L85").Value = 28660
$ws.Range("M85").Value = -3563.25
$ws.Range("N85").Value = -31312

$ws.Range("H99").Value = 2150.5557
$ws.Range("I99").Value = 1438.3334
$ws.Range("J99").Value = 3575
$ws.Range("K99").Value = 1438.3334
$ws.Range("L99").Value = 3575
$ws.Range("M99").Value = 59.66660000000002
$ws.Range("N99").Value = -6571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 909.2381
$ws.Range("I5").Value = 815.6923
$ws.Range("J5").Value = 1061.25
$ws.Range("K5").Value = 2447.0769
$ws.Range("L5").Value = 3183.75
$ws.Range("M5").Value = -2335.0769
$ws.Range("N5").Value = -3407.75

$ws.Range("H122").Value = 517.5417
$ws.Range("I122").Value = 326.92307
$ws.Range("J122").Value = 742.8182
$ws.Range("K122").Value = 2942.30763
$ws.Range("L122").Value = 6685.3638
$ws.Range("M122").Value = -492.3076299999998
$ws.Range("N122").Value = -11585.3638

$ws.Range("H131").Value = 17858068
$ws.Range("J131").Value = 26316676
$ws.Range("L131").Value = 78950028
$ws.Range("N131").Value = -78960108

$ws.Range("H132").Value = 1458.9062
$ws.Range("I132").Value = 950.4286
$ws.Range("J132").Value = 2429.6365
$ws.Range("K132").Value = 8553.857399999999
$ws.Range("L132").Value = 21866.7285
$ws.Range("M132").Value = -6023.857399999999
$ws.Range("N132").Value = -26926.7285

$ws.Range("H135").Value = 909.2381
$ws.Range("I135").Value = 815.6923
$ws.Range("J135").Value = 1061.25
$ws.Range("K135").Value = 7341.2307
$ws.Range("L135").Value = 9551.25
$ws.Range("M135").Value = -4806.2307
$ws.Range("N135").Value = -14621.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1909.4445
$ws.Range("I97").Value = 1396.25
$ws.Range("J97").Value = 2320
$ws.Range("K97").Value = 1396.25
$ws.Range("L97").Value = 2320
$ws.Range("M97").Value = -900.25
$ws.Range("N97").Value = -3312

$ws.Range("H102").Value = 1566.619
$ws.Range("I102").Value = 1025.8334
$ws.Range("J102").Value = 2287.6667
$ws.Range("K102").Value = 1025.8334
$ws.Range("L102").Value = 2287.6667
$ws.Range("M102").Value = 596.1666
$ws.Range("N102").Value = -5531.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2434
$ws.Range("I100").Value = 1389.5555
$ws.Range("K100").Value = 1389.5555
$ws.Range("M100").Value = -848.5554999999999

$ws.Range("H122").Value = 2324.1052
$ws.Range("I122").Value = 2432.7144
$ws.Range("J122").Value = 2020
$ws.Range("K122").Value = 7298.1432
$ws.Range("L122").Value = 6060
$ws.Range("M122").Value = -4848.1432
$ws.Range("N122").Value = -10960

$ws.Range("H136").Value = 40007108
$ws.Range("I136").Value = 8093.2
$ws.Range("J136").Value = 100005630
$ws.Range("K136").Value = 24279.6
$ws.Range("L136").Value = 300016890
$ws.Range("M136").Value = -21729.6
$ws.Range("N136").Value = -300021990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
